$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("D3").Value = "2016-01-11 05:02:20"
$ws.Range("G3").Value = "2016-01-11 05:03:27"

$ws2 = $wb.Worksheets.Item("de-de")
$ws2.Range("D3").Value = "2016-01-11 05:02:37"
$ws2.Range("G3").Value = "2016-01-11 05:03:54"
